$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(" MVC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$templateFT = $rng.FormattedText
$destPos = $rng.End

$endRng = $d.Content
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
$scratchStart = $d.Content.End - 1
$scratch = $d.Range($scratchStart, $scratchStart)
$scratch.FormattedText = $templateFT
$scratch2 = $d.Range($scratchStart, $scratchStart + 4)
$scratch2.Text = " and bootstrap"
$newFT = $scratch2.FormattedText

$dest = $d.Range($destPos, $destPos)
$dest.FormattedText = $newFT

# cleanup: delete scratch paragraph #74 in its entirety (incl. paragraph mark)
Write-Output "ParaCount before cleanup=$($d.Paragraphs.Count)"
$scratchPara = $d.Paragraphs.Item($d.Paragraphs.Count)
Write-Output "scratchPara Start=$($scratchPara.Range.Start) End=$($scratchPara.Range.End) Text=[$($scratchPara.Range.Text)]"
$scratchPara.Range.Delete()
Write-Output "ParaCount after cleanup=$($d.Paragraphs.Count)"
